# A new price-quote row was recorded for "Agrícola del Norte S.A. de Arica"
# (Cardinal / 1a (cosecha), Región de Coquimbo) that belongs right above the
# existing row 63 entry (chronologically/by sort order within the sheet).
# Insert a fresh row at position 63 - this shifts every following row
# (old 63..111) down by one (new 64..112), exactly like typing a new line in
# the middle of the data table - and then fill in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(63).Insert()

$ws.Cells.Item(63, 1).Value  = 1
$ws.Cells.Item(63, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(63, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(63, 4).Value  = 44907
$ws.Cells.Item(63, 5).Value  = 15
$ws.Cells.Item(63, 6).Value  = 100114001
$ws.Cells.Item(63, 7).Value  = "Papa"
$ws.Cells.Item(63, 8).Value  = "Cardinal"
$ws.Cells.Item(63, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(63, 10).Value = 1150
$ws.Cells.Item(63, 11).Value = 18000
$ws.Cells.Item(63, 12).Value = 19000
$ws.Cells.Item(63, 13).Value = 18565
$ws.Cells.Item(63, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(63, 16).Value = 743
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
